$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: task "ADD MORE IMAGE PROCESSING" - completed.
# Set the completion date (TANGGAL SELESAI, column D) for row 9.
$ws.Range("D9").Value = "7 Desember 2021"

# Update status (column F) from "ON GOING" to "DONE".
$ws.Range("F9").Value = "DONE"

# Move the active selection to D12, matching the saved cursor position.
$ws.Range("D12").Select()
